{"js": "// Update the date line (first paragraph of the body).\nconst firstPara = context.document.body.paragraphs.getFirst();\nfirstPara.insertText(\"2025-01-24 Friday\", Word.InsertLocation.replace);\n\n// Update each of the 100 table cells (20 rows x 5 columns) with new\n// expressions, in row-major order, matching the original cell order exactly.\nconst newValues = [\n  \"27+41=\", \"98-66=\", \"85-67=\", \"66-13=\", \"55-5=\",\n  \"33-29=\", \"32+27=\", \"82-67=\", \"44-34=\", \"79-63=\",\n  \"50-1=\", \"62-37=\", \"82-67=\", \"45+26=\", \"96-57=\",\n  \"36-5=\", \"50+39=\", \"98-55=\", \"93+1=\", \"73-2=\",\n  \"57+17=\", \"90-22=\", \"42+27=\", \"12+44=\", \"91-42=\",\n  \"10+9=\", \"14+10=\", \"37-20=\", \"54-30=\", \"99-61=\",\n  \"2+62=\", \"60-34=\", \"61+13=\", \"8+11=\", \"51-22=\",\n  \"25+37=\", \"9+42=\", \"65-27=\", \"97-75=\", \"21+21=\",\n  \"3+34=\", \"24+58=\", \"78+12=\", \"93-78=\", \"85-69=\",\n  \"7+61=\", \"44-5=\", \"92-52=\", \"40-28=\", \"32+9=\",\n  \"71-9=\", \"56+41=\", \"60-9=\", \"62-38=\", \"0+12=\",\n  \"52-8=\", \"51+2=\", \"38+60=\", \"54-49=\", \"27+51=\",\n  \"12-0=\", \"50-16=\", \"10+33=\", \"94-85=\", \"42+18=\",\n  \"84-38=\", \"68-55=\", \"26+22=\", \"9+53=\", \"23-13=\",\n  \"95-57=\", \"54-7=\", \"28+39=\", \"20+46=\", \"76-61=\",\n  \"42+43=\", \"9+79=\", \"45+44=\", \"76-1=\", \"82-20=\",\n  \"27+19=\", \"81-4=\", \"35-22=\", \"63+27=\", \"11+38=\",\n  \"4+59=\", \"69+24=\", \"98-42=\", \"44+30=\", \"78-51=\",\n  \"37-27=\", \"90-44=\", \"14+84=\", \"82-17=\", \"9+41=\",\n  \"63-33=\", \"5-0=\", \"29+35=\", \"40+26=\", \"69+0=\"\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst colCount = 5;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r * colCount + c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = \"2025-01-24 Friday\"\n\n# Update each of the 100 table cells (20 rows x 5 columns) with new expressions,\n# in row-major order, matching the original cell order exactly.\n$newValues = @(\n    \"27+41=\",\n    \"98-66=\",\n    \"85-67=\",\n    \"66-13=\",\n    \"55-5=\",\n    \"33-29=\",\n    \"32+27=\",\n    \"82-67=\",\n    \"44-34=\",\n    \"79-63=\",\n    \"50-1=\",\n    \"62-37=\",\n    \"82-67=\",\n    \"45+26=\",\n    \"96-57=\",\n    \"36-5=\",\n    \"50+39=\",\n    \"98-55=\",\n    \"93+1=\",\n    \"73-2=\",\n    \"57+17=\",\n    \"90-22=\",\n    \"42+27=\",\n    \"12+44=\",\n    \"91-42=\",\n    \"10+9=\",\n    \"14+10=\",\n    \"37-20=\",\n    \"54-30=\",\n    \"99-61=\",\n    \"2+62=\",\n    \"60-34=\",\n    \"61+13=\",\n    \"8+11=\",\n    \"51-22=\",\n    \"25+37=\",\n    \"9+42=\",\n    \"65-27=\",\n    \"97-75=\",\n    \"21+21=\",\n    \"3+34=\",\n    \"24+58=\",\n    \"78+12=\",\n    \"93-78=\",\n    \"85-69=\",\n    \"7+61=\",\n    \"44-5=\",\n    \"92-52=\",\n    \"40-28=\",\n    \"32+9=\",\n    \"71-9=\",\n    \"56+41=\",\n    \"60-9=\",\n    \"62-38=\",\n    \"0+12=\",\n    \"52-8=\",\n    \"51+2=\",\n    \"38+60=\",\n    \"54-49=\",\n    \"27+51=\",\n    \"12-0=\",\n    \"50-16=\",\n    \"10+33=\",\n    \"94-85=\",\n    \"42+18=\",\n    \"84-38=\",\n    \"68-55=\",\n    \"26+22=\",\n    \"9+53=\",\n    \"23-13=\",\n    \"95-57=\",\n    \"54-7=\",\n    \"28+39=\",\n    \"20+46=\",\n    \"76-61=\",\n    \"42+43=\",\n    \"9+79=\",\n    \"45+44=\",\n    \"76-1=\",\n    \"82-20=\",\n    \"27+19=\",\n    \"81-4=\",\n    \"35-22=\",\n    \"63+27=\",\n    \"11+38=\",\n    \"4+59=\",\n    \"69+24=\",\n    \"98-42=\",\n    \"44+30=\",\n    \"78-51=\",\n    \"37-27=\",\n    \"90-44=\",\n    \"14+84=\",\n    \"82-17=\",\n    \"9+41=\",\n    \"63-33=\",\n    \"5-0=\",\n    \"29+35=\",\n    \"40+26=\",\n    \"69+0=\"\n)\n\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx = $idx + 1\n    }\n}\n"}
